$wb = $excel.ActiveWorkbook

# --- Beads sheet: update control file paths to include sub-folder per control type ---
$wsBeads = $wb.Worksheets.Item("Beads")
$wsBeads.Range("C3").Value = "FCFiles/controls/sfc2/sample001.fcs"
$wsBeads.Range("C4").Value = "FCFiles/controls/sfc1/sample003.fcs"
$wsBeads.Range("C5").Value = "FCFiles/controls/nfc/sample003.fcs"

# --- Samples sheet: update control rows (NFC / SFC1 / SFC2) ---
$wsSamples = $wb.Worksheets.Item("Samples")

# NFC row (row 2): file path now lives under controls/nfc/
$wsSamples.Range("D2").Value = "FCFiles/controls/nfc/sample004.fcs"

# SFC1 row (row 3): new strain (sJS1061), new plasmid description, new file path
$wsSamples.Range("D3").Value = "FCFiles/controls/sfc1/sample007.fcs"
$wsSamples.Range("H3").Value = "sJS1061"
$wsSamples.Range("J3").Value = "pSC31_3, pJS0306"

# SFC2 row (row 4): file path now lives under controls/sfc2/
$wsSamples.Range("D4").Value = "FCFiles/controls/sfc2/sample019.fcs"

# --- Restore selections to match the saved view state ---
# Select a cell on the Beads sheet first (this will become the non-active tab's
# remembered selection), then finish on the Samples sheet so it stays the active tab.
$wsBeads.Range("C5").Select()
$wsSamples.Range("I24").Select()
